# Auto-generated PowerShell COM-interop script
# Applies numeric cell value updates to the Spriggan Profits workbook sheets
# (ALC, ARM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 1533.3334
$ws.Range("I18").Value = 1300
$ws.Range("K18").Value = 1300
$ws.Range("M18").Value = -1016
# Row 40
$ws.Range("H40").Value = 6176510
$ws.Range("I40").Value = 3886.5293
$ws.Range("J40").Value = 111111110
$ws.Range("K40").Value = 3886.5293
$ws.Range("L40").Value = 111111110
$ws.Range("M40").Value = -3711.5293
$ws.Range("N40").Value = -111111460
# Row 41
$ws.Range("H41").Value = 1430.1538
$ws.Range("I41").Value = 1681.9
$ws.Range("K41").Value = 1681.9
$ws.Range("M41").Value = -1241.9
# Row 70
$ws.Range("H70").Value = 15066.667
$ws.Range("I70").Value = 3233.3333
$ws.Range("J70").Value = 26900
$ws.Range("K70").Value = 9699.999899999999
$ws.Range("L70").Value = 80700
$ws.Range("M70").Value = -9429.999899999999
$ws.Range("N70").Value = -81240
# Row 73
$ws.Range("H73").Value = 15066.667
$ws.Range("I73").Value = 3233.3333
$ws.Range("J73").Value = 26900
$ws.Range("K73").Value = 9699.999899999999
$ws.Range("L73").Value = 80700
$ws.Range("M73").Value = -8763.999899999999
$ws.Range("N73").Value = -82572
# Row 98
$ws.Range("H98").Value = 1565
$ws.Range("I98").Value = 1598.8235
$ws.Range("K98").Value = 1598.8235
$ws.Range("M98").Value = -100.8235
# Row 122
$ws.Range("H122").Value = 1565
$ws.Range("I122").Value = 1598.8235
$ws.Range("K122").Value = 4796.470499999999
$ws.Range("M122").Value = -2346.470499999999
# Row 136
$ws.Range("H136").Value = 50999
$ws.Range("J136").Value = 50999
$ws.Range("L136").Value = 50999
$ws.Range("N136").Value = -61199
# Row 140
$ws.Range("H140").Value = 54982.5
$ws.Range("J140").Value = 54982.5
$ws.Range("L140").Value = 54982.5
$ws.Range("N140").Value = -65342.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 33746.332
$ws.Range("J4").Value = 33746.332
$ws.Range("L4").Value = 33746.332
$ws.Range("N4").Value = -33978.332
# Row 45
$ws.Range("H45").Value = 972.75
$ws.Range("I45").Value = 972.75
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 972.75
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -595.75
$ws.Range("N45").ClearContents()
# Row 61
$ws.Range("H61").Value = 23260476
$ws.Range("I61").Value = 25004490
$ws.Range("K61").Value = 25004490
$ws.Range("M61").Value = -25004278
# Row 110
$ws.Range("H110").Value = 58392.11
$ws.Range("I110").Value = 79081.30499999999
$ws.Range("J110").Value = 4600.2
$ws.Range("K110").Value = 79081.30499999999
$ws.Range("L110").Value = 4600.2
$ws.Range("M110").Value = -77036.30499999999
$ws.Range("N110").Value = -8690.200000000001
# Row 136
$ws.Range("H136").Value = 23260476
$ws.Range("I136").Value = 25004490
$ws.Range("K136").Value = 75013470
$ws.Range("M136").Value = -75010920

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8108.4
$ws.Range("I31").Value = 5524.212
$ws.Range("J31").Value = 11984.682
$ws.Range("K31").Value = 5524.212
$ws.Range("L31").Value = 11984.682
$ws.Range("M31").Value = -5229.212
$ws.Range("N31").Value = -12574.682
# Row 34
$ws.Range("H34").Value = 8108.4
$ws.Range("I34").Value = 5524.212
$ws.Range("J34").Value = 11984.682
$ws.Range("K34").Value = 5524.212
$ws.Range("L34").Value = 11984.682
$ws.Range("M34").Value = -5322.212
$ws.Range("N34").Value = -12388.682
# Row 141
$ws.Range("H141").Value = 565000
$ws.Range("J141").Value = 720000
$ws.Range("L141").Value = 720000
$ws.Range("N141").Value = -730360

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 15
$ws.Range("H15").Value = 135
$ws.Range("I15").Value = 113.333336
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 340.000008
$ws.Range("L15").Value = 600
$ws.Range("M15").Value = -200.000008
$ws.Range("N15").Value = -880
# Row 44
$ws.Range("H44").Value = 147.85715
$ws.Range("J44").Value = 325
$ws.Range("L44").Value = 975
$ws.Range("N44").Value = -1771
# Row 130
$ws.Range("H130").Value = 5071.2856
$ws.Range("I130").Value = 2000
$ws.Range("J130").Value = 5583.1665
$ws.Range("K130").Value = 6000
$ws.Range("L130").Value = 16749.4995
$ws.Range("M130").Value = -980
$ws.Range("N130").Value = -26789.4995
# Row 131
$ws.Range("H131").Value = 1631.381
$ws.Range("J131").Value = 1836.25
$ws.Range("L131").Value = 5508.75
$ws.Range("N131").Value = -15588.75

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 100000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 100000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -100224
# Row 70
$ws.Range("H70").Value = 10301.5
$ws.Range("I70").Value = 9646.6
$ws.Range("K70").Value = 9646.6
$ws.Range("M70").Value = -9376.6
# Row 73
$ws.Range("H73").Value = 10301.5
$ws.Range("I73").Value = 9646.6
$ws.Range("K73").Value = 9646.6
$ws.Range("M73").Value = -8710.6
# Row 80
$ws.Range("H80").Value = 2513.7144
$ws.Range("I80").Value = 2372
$ws.Range("K80").Value = 2372
$ws.Range("M80").Value = -1374
# Row 83
$ws.Range("H83").Value = 2513.7144
$ws.Range("I83").Value = 2372
$ws.Range("K83").Value = 11860
$ws.Range("M83").Value = -6868
# Row 132
$ws.Range("H132").Value = 2406740.2
$ws.Range("I132").Value = 2661864
$ws.Range("K132").Value = 7985592
$ws.Range("M132").Value = -7983062

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 19
$ws.Range("H19").Value = 537
$ws.Range("J19").Value = 600
$ws.Range("L19").Value = 600
$ws.Range("N19").Value = -940
# Row 100
$ws.Range("H100").Value = 9973778
$ws.Range("I100").Value = 16617797
$ws.Range("J100").Value = 7749.875
$ws.Range("K100").Value = 16617797
$ws.Range("L100").Value = 7749.875
$ws.Range("M100").Value = -16617256
$ws.Range("N100").Value = -8831.875
# Row 136
$ws.Range("H136").Value = 2921.2222
$ws.Range("I136").Value = 1099
$ws.Range("J136").Value = 3028.4119
$ws.Range("K136").Value = 3297
$ws.Range("L136").Value = 9085.235700000001
$ws.Range("M136").Value = -747
$ws.Range("N136").Value = -14185.2357

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 11365291
$ws.Range("I132").Value = 14287152
$ws.Range("K132").Value = 42861456
$ws.Range("M132").Value = -42858926
# Row 136
$ws.Range("H136").Value = 27778614
$ws.Range("J136").Value = 1415.8
$ws.Range("L136").Value = 4247.4
$ws.Range("N136").Value = -9347.4
